$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels to lowercase (name/session/staff_id)
$ws.Range("C1").Value = "session"
$ws.Range("D1").Value = "staff_id"
$ws.Range("B1").Value = "name"

# Bump the ID column (A2:A9) by 1 so it runs 1..8 instead of 0..7
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value()
    $cell.Value = $old + 1
}

# Add the generated INSERT statement formulas in column E
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=CONCAT(""INSERT INTO subject (name,session,staff_id) VALUE ('"",B$r,""','"",C$r,""','"",D$r,""');"")"
}

# Autofit column E like the bestFit width recorded in the diff
$ws.Columns.Item(5).AutoFit() | Out-Null

$ws.Range("B2").Select() | Out-Null
